$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen the Stimulus column (C) to fit the new, longer sentence; the
# Response column (D) no longer needs to be wide since it is now blank
# for the new TR control-code rows.
$ws.Columns("C").ColumnWidth = 62.44140625
$ws.Columns("D").ColumnWidth = 8.77734375

# Row 3: new TR control-code stimulus replaces the old one, and the
# associated response/index/reaction-time are cleared (no response was
# collected for this new control trial).
$ws.Range("C3").Value2 = "O queijo azul tem um sabor intenso e peculiar, ao contrário do queijo verde."
$ws.Range("D3:F3").Value2 = ""
$ws.Range("D3:F3").Style = "Normal"

# Row 4: stimulus text stays the same, but the recorded response/index/
# reaction-time event data is cleared out as well.
$ws.Range("D4:F4").Value2 = ""
$ws.Range("D4:F4").Style = "Normal"

# Row 5: this trial is reset back to an empty/unused event row, matching
# the blank template rows further down the log.
$ws.Range("A5").Value2 = 0
$ws.Range("B5").Value2 = 0
$ws.Range("C5").Value2 = ""
$ws.Range("D5").Value2 = ""
$ws.Range("E5").Value2 = 0
$ws.Range("F5").Value2 = 0
$ws.Range("G5").Value2 = ""
$ws.Range("C5:D5").Style = "Normal"
$ws.Range("G5").Style = "Normal"
